$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific cells to Text format so numeric-looking strings
# (e.g. "554.16", "0.999") are preserved exactly as text, matching
# the source inline-string cell type instead of being parsed as numbers.
$textCells = @("D5","D6","D7","D10","D12","D14","D16","D19","D20","D21","D22","D23","D24","D30","D31","D32","D33","D35","D37","D38","D39","D40","D42","D43","D44","D46","D47","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet order.
$ws.Range("D2").Value = "61.881.00"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "2.584.71"
$ws.Range("E3").Value = "  -4.60%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "554.16"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "154.03"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").Value = "0.364"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").Value = "3.038.53"
$ws.Range("E13").Value = "  -4.80%  "
$ws.Range("D14").Value = "25.47"
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").Value = "61.793.72"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").Value = "0.0000144"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "2.582.36"
$ws.Range("E17").Value = "  -4.97%  "
$ws.Range("E18").Value = "  -4.23%  "
$ws.Range("D19").Value = "4.53"
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").Value = "338.70"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  -6.07%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "0.495"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("D24").Value = "62.82"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "0.0₃0839"
$ws.Range("E28").Value = "  -4.62%  "
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "7.06"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").Value = "1.30"
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("D32").Value = "160.40"
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("D35").Value = "4.70"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("E36").Value = "  -3.68%  "
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").Value = "338.10"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("D39").Value = "6.07"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "0.895"
$ws.Range("E40").Value = "  -5.75%  "
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("D42").Value = "37.49"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "20.57"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "2.136.24"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("D46").Value = "0.608"
$ws.Range("E46").Value = "  -3.94%  "
$ws.Range("D47").Value = "10.93"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "19.67"
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0546"
$ws.Range("E49").Value = "  -4.18%  "
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  -1.50%  "
